$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2026-02-24 20:48:21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2026-02-24 20:48:24"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "32%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2026-02-24 20:48:26"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "71%"
$ws.Range("O4").NumberFormat = "@"
$ws.Range("O4").Value = "13.0 °C"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2026-02-24 20:48:29"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2026-02-24 20:48:31"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "69%"
$ws.Range("O6").NumberFormat = "@"
$ws.Range("O6").Value = "14.1 °C"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2026-02-24 20:48:34"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "71%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2026-02-24 20:48:36"
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "1019.6 hPa"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2026-02-24 20:48:39"
$ws.Range("O9").NumberFormat = "@"
$ws.Range("O9").Value = "11.9 °C"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2026-02-24 20:48:41"
$ws.Range("O10").NumberFormat = "@"
$ws.Range("O10").Value = "11.3 °C"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2026-02-24 20:48:44"
$ws.Range("O11").NumberFormat = "@"
$ws.Range("O11").Value = "9.2 °C"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2026-02-24 20:48:46"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "91%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2026-02-24 20:48:48"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2026-02-24 20:48:51"
$ws.Range("N14").NumberFormat = "@"
$ws.Range("N14").Value = "7.5 °C 20:29 TU"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2026-02-24 20:48:53"
$ws.Range("O15").NumberFormat = "@"
$ws.Range("O15").Value = "12.0 °C"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2026-02-24 20:48:56"
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = "25.6 km/h - 195º 20:20 TU"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2026-02-24 20:48:58"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2026-02-24 20:49:01"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2026-02-24 20:49:03"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2026-02-24 20:49:06"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "35%"
$ws.Range("K20").NumberFormat = "@"
$ws.Range("K20").Value = "15.0 MJ/m2"
$ws.Range("O20").NumberFormat = "@"
$ws.Range("O20").Value = "3.6 °C"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2026-02-24 20:49:08"
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value = "1021.9 hPa"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2026-02-24 20:49:11"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "25%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2026-02-24 20:49:13"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "23%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2026-02-24 20:49:16"
$ws.Range("J24").NumberFormat = "@"
$ws.Range("J24").Value = "1021.2 hPa"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2026-02-24 20:49:18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2026-02-24 20:49:21"
$ws.Range("O26").NumberFormat = "@"
$ws.Range("O26").Value = "11.8 °C"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2026-02-24 20:49:23"
$ws.Range("K27").NumberFormat = "@"
$ws.Range("K27").Value = "15.2 MJ/m2"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "2026-02-24 20:49:26"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "2026-02-24 20:49:28"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "2026-02-24 20:49:31"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "75%"
$ws.Range("J30").NumberFormat = "@"
$ws.Range("J30").Value = "1019.7 hPa"
$ws.Range("O30").NumberFormat = "@"
$ws.Range("O30").Value = "13.2 °C"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "2026-02-24 20:49:33"
$ws.Range("J31").NumberFormat = "@"
$ws.Range("J31").Value = "1019.1 hPa"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "2026-02-24 20:49:36"
$ws.Range("O32").NumberFormat = "@"
$ws.Range("O32").Value = "7.3 °C"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "2026-02-24 20:49:38"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "2026-02-24 20:49:41"
$ws.Range("O34").NumberFormat = "@"
$ws.Range("O34").Value = "4.8 °C"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "2026-02-24 20:49:43"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "36%"
$ws.Range("J35").NumberFormat = "@"
$ws.Range("J35").Value = "1020.4 hPa"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "2026-02-24 20:49:46"
$ws.Range("O36").NumberFormat = "@"
$ws.Range("O36").Value = "13.0 °C"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "2026-02-24 20:49:48"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "71%"
$ws.Range("O37").NumberFormat = "@"
$ws.Range("O37").Value = "8.8 °C"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "2026-02-24 20:49:51"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "72%"
$ws.Range("O38").NumberFormat = "@"
$ws.Range("O38").Value = "12.1 °C"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2026-02-24 20:49:53"
$ws.Range("N39").NumberFormat = "@"
$ws.Range("N39").Value = "1.3 °C 20:15 TU"
$ws.Range("O39").NumberFormat = "@"
$ws.Range("O39").Value = "4.5 °C"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2026-02-24 20:49:56"
$ws.Range("J40").NumberFormat = "@"
$ws.Range("J40").Value = "1022.6 hPa"
$ws.Range("O40").NumberFormat = "@"
$ws.Range("O40").Value = "8.7 °C"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2026-02-24 20:49:58"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "79%"
$ws.Range("J41").NumberFormat = "@"
$ws.Range("J41").Value = "1020.5 hPa"
$ws.Range("O41").NumberFormat = "@"
$ws.Range("O41").Value = "10.7 °C"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2026-02-24 20:50:01"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2026-02-24 20:50:03"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2026-02-24 20:50:05"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2026-02-24 20:50:08"
$ws.Range("L45").NumberFormat = "@"
$ws.Range("L45").Value = "26.3 km/h - 120º 20:04 TU"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2026-02-24 20:50:10"
$ws.Range("J46").NumberFormat = "@"
$ws.Range("J46").Value = "1021.2 hPa"
$ws.Range("O46").NumberFormat = "@"
$ws.Range("O46").Value = "10.7 °C"
